# Append two new rows of scraped data (rows 3 and 4) to the BIIBNamed sheet,
# and widen column A slightly to keep the "bestFit" date column readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 42600.835023148145
$ws.Range("B3").Value = "Named"
$ws.Range("C3").Value = 11022
$ws.Range("D3").Value = 4882
$ws.Range("E3").Value = 319
$ws.Range("F3").Value = 51
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = 78
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

# Row 4
$ws.Range("A4").Value = 42600.88
$ws.Range("B4").Value = "Named"
$ws.Range("C4").Value = 10141
$ws.Range("D4").Value = 5097
$ws.Range("E4").Value = 334
$ws.Range("F4").Value = 55
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 78
$ws.Range("I4").Value = 21
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0

# Widen column A (date column) slightly to keep the best-fit content visible
$ws.Columns.Item(1).ColumnWidth = 14
